# Daily COVID stats update "po 22. 03. 2021" (as of 22 Mar 2021):
#  - revises the AgTests (col F) and AgPosit (col G) figures for a number of
#    already-present days (rows 334-379), and
#  - appends two brand-new daily rows (380, 381) for 2021-03-19 and 2021-03-20.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Revisions to existing rows (only columns F and G change) ---------------
$updates = @(
    @{Row=334; F=196012; G=3471}
    @{Row=335; F=131099; G=3009}
    @{Row=337; F=104039; G=$null}
    @{Row=338; F=227436; G=3189}
    @{Row=339; F=658497; G=$null}
    @{Row=341; F=291244; G=$null}
    @{Row=343; F=132674; G=2971}
    @{Row=344; F=135698; G=2488}
    @{Row=345; F=291801; G=3315}
    @{Row=347; F=342374; G=2907}
    @{Row=348; F=232253; G=3250}
    @{Row=349; F=159571; G=2758}
    @{Row=352; F=306779; G=3536}
    @{Row=353; F=721189; G=5259}
    @{Row=354; F=310395; G=2846}
    @{Row=355; F=222200; G=3456}
    @{Row=356; F=160079; G=2878}
    @{Row=357; F=138579; G=$null}
    @{Row=358; F=157339; G=2599}
    @{Row=362; F=228066; G=3170}
    @{Row=363; F=188072; G=2759}
    @{Row=364; F=167925; G=2474}
    @{Row=365; F=180165; G=2352}
    @{Row=366; F=337206; G=2831}
    @{Row=367; F=760491; G=3888}
    @{Row=369; F=232983; G=2570}
    @{Row=370; F=181876; G=2031}
    @{Row=371; F=158308; G=1939}
    @{Row=372; F=175810; G=1830}
    @{Row=373; F=343576; G=2344}
    @{Row=374; F=761912; G=3359}
    @{Row=375; F=349788; G=1849}
    @{Row=376; F=218381; G=2186}
    @{Row=377; F=174477; G=1814}
    @{Row=378; F=153018; G=1497}
    @{Row=379; F=170552; G=1553}
)

foreach ($u in $updates) {
    $r = $u.Row
    $ws.Cells.Item($r, 6).Value = $u.F
    if ($null -ne $u.G) {
        $ws.Cells.Item($r, 7).Value = $u.G
    }
}

# --- Brand-new rows appended at the bottom of the table ----------------------
$newRows = @(
    @{Row=380; A=44274; B=347944; C=13293; D=1795; E=8978; F=318464; G=1830}
    @{Row=381; A=44275; B=348869; C=5896;  D=925;  E=9044; F=620497; G=2646}
)

foreach ($nr in $newRows) {
    $r = $nr.Row
    $ws.Cells.Item($r, 1).Value = $nr.A
    $ws.Cells.Item($r, 2).Value = $nr.B
    $ws.Cells.Item($r, 3).Value = $nr.C
    $ws.Cells.Item($r, 4).Value = $nr.D
    $ws.Cells.Item($r, 5).Value = $nr.E
    $ws.Cells.Item($r, 6).Value = $nr.F
    $ws.Cells.Item($r, 7).Value = $nr.G
}

Write-Output "Updated rows 334-379 and appended rows 380-381"
